$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '30.476.07'
Set-TextValue 2 5 '  +0.10%  '

# Row 3
Set-TextValue 3 4 '1.933.92'
Set-TextValue 3 5 '  +4.38%  '

# Row 4
Set-TextValue 4 4 '0.9990'
Set-TextValue 4 5 '  -0.10%  '

# Row 5
Set-TextValue 5 4 '240.62'
Set-TextValue 5 5 '  +3.06%  '

# Row 6
Set-TextValue 6 4 '0.9996'
Set-TextValue 6 5 '  -0.05%  '

# Row 7
Set-TextValue 7 4 '0.4744'
Set-TextValue 7 5 '  -0.16%  '

# Row 8
Set-TextValue 8 4 '0.2872'
Set-TextValue 8 5 '  +4.32%  '

# Row 9
Set-TextValue 9 4 '0.06614'
Set-TextValue 9 5 '  +4.28%  '

# Row 10
Set-TextValue 10 4 '19.25'
Set-TextValue 10 5 '  +7.99%  '

# Row 11
Set-TextValue 11 4 '107.18'
Set-TextValue 11 5 '  +26.46%  '

# Row 12
Set-TextValue 12 4 '1.922.53'
Set-TextValue 12 5 '  -3.16%  '

# Row 13
Set-TextValue 13 4 '0.07632'
Set-TextValue 13 5 '  +2.29%  '

# Row 14
Set-TextValue 14 4 '5.164'
Set-TextValue 14 5 '  +4.11%  '

# Row 15
Set-TextValue 15 4 '0.6641'
Set-TextValue 15 5 '  +6.28%  '

# Row 16
Set-TextValue 16 4 '308.59'
Set-TextValue 16 5 '  +25.45%  '

# Row 17
Set-TextValue 17 4 '30.489.49'
Set-TextValue 17 5 '  +0.28%  '

# Row 18
Set-TextValue 18 4 '13.05'
Set-TextValue 18 5 '  +2.92%  '

# Row 19
Set-TextValue 19 4 '0.9993'
Set-TextValue 19 5 '  -0.07%  '

# Row 20
Set-TextValue 20 4 '0.000007562'
Set-TextValue 20 5 '  +3.19%  '

# Row 21
Set-TextValue 21 4 '2.181.34'
Set-TextValue 21 5 '  +3.62%  '

# Row 22
Set-TextValue 22 4 '5.318'
Set-TextValue 22 5 '  +8.11%  '

# Row 23
Set-TextValue 23 4 '0.9991'
Set-TextValue 23 5 '  -0.17%  '

# Row 24
Set-TextValue 24 4 '6.304'
Set-TextValue 24 5 '  +6.53%  '

# Row 25
Set-TextValue 25 4 '9.304'
Set-TextValue 25 5 '  +2.61%  '

# Row 26
Set-TextValue 26 4 '167.61'
Set-TextValue 26 5 '  +1.96%  '

# Row 27
Set-TextValue 27 4 '20.36'
Set-TextValue 27 5 '  +13.28%  '

# Row 28
Set-TextValue 28 4 '2.060'
Set-TextValue 28 5 '  +9.80%  '

# Row 29
Set-TextValue 29 4 '0.1111'
Set-TextValue 29 5 '  +8.15%  '

# Row 30
Set-TextValue 30 4 '1.366'
Set-TextValue 30 5 '  +1.37%  '

# Row 31
Set-TextValue 31 4 '4.118'
Set-TextValue 31 5 '  +1.90%  '

# Row 32
Set-TextValue 32 5 '  +2.88%  '

# Row 33
Set-TextValue 33 4 '0.05046'
Set-TextValue 33 5 '  +4.41%  '

# Row 34
Set-TextValue 34 4 '0.7435'
Set-TextValue 34 5 '  +6.40%  '

# Row 35
Set-TextValue 35 4 '1.155'
Set-TextValue 35 5 '  +2.30%  '

# Row 36
Set-TextValue 36 4 '2.746'
Set-TextValue 36 5 '  +1.53%  '

# Row 37
Set-TextValue 37 4 '0.01970'
Set-TextValue 37 5 '  +3.37%  '

# Row 38
Set-TextValue 38 4 '2.689'
Set-TextValue 38 5 '  +0.25%  '

# Row 39
Set-TextValue 39 4 '2.055'
Set-TextValue 39 5 '  +3.30%  '

# Row 40
Set-TextValue 40 4 '0.8799'
Set-TextValue 40 5 '  +0.25%  '

# Row 41
Set-TextValue 41 4 '107.42'
Set-TextValue 41 5 '  +0.69%  '

# Row 42
Set-TextValue 42 4 '70.86'
Set-TextValue 42 5 '  +11.90%  '

# Row 43
Set-TextValue 43 4 '5.807'
Set-TextValue 43 5 '  +5.62%  '

# Row 44
Set-TextValue 44 4 '0.9994'
Set-TextValue 44 5 '  -0.06%  '

# Row 45
Set-TextValue 45 4 '0.4189'
Set-TextValue 45 5 '  +3.28%  '

# Row 46
Set-TextValue 46 4 '7.304'
Set-TextValue 46 5 '  +1.89%  '

# Row 47
Set-TextValue 47 4 '9.281'
Set-TextValue 47 5 '  +7.98%  '

# Row 48
Set-TextValue 48 2 'Algorand'
Set-TextValue 48 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 48 4 '0.1212'
Set-TextValue 48 5 '  +1.00%  '

# Row 49
Set-TextValue 49 2 'Elrond'
Set-TextValue 49 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 49 4 '34.88'
Set-TextValue 49 5 '  +2.62%  '

# Row 50
Set-TextValue 50 4 '0.05622'
Set-TextValue 50 5 '  +2.25%  '

# Row 51
Set-TextValue 51 4 '0.3857'
Set-TextValue 51 5 '  +4.66%  '
